$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 10)
$ws.Range("A10").Value = 44509
$ws.Range("B10").Value = 2.36
$ws.Range("C10").Value = "UI Bug fixen und Media Kontrollen Notification"

# Apply a thin left/right border to B10 (new border/style combination)
$rng = $ws.Range("B10")
$rng.Borders.Item(7).LineStyle = 1
$rng.Borders.Item(10).LineStyle = 1

# Update the selected cell to C10
$ws.Range("C10").Select()

# Update the saved window position
$excel.ActiveWindow.Left = 3765
$excel.ActiveWindow.Top = 4215
